$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("courses")

# Update the department name (C2) to the shorter label
$ws.Range("C2").Value = "Community Services"

# Clear the promotion validity text (R2) but keep its formatting
$ws.Range("R2").ClearContents()
